# Update workbook to add carjacking data for 2022-12-05 (extends the
# "through November 26" snapshot to "through November 27").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet and update the header/title string -------------------
$ws.Name = "Through 2022-11-27"
$ws.Range("B1").Value = "November 2022 (through November 27)"

# --- Per-neighborhood cell updates (month columns for Nov 2022 / Nov 2021 /
#     Nov 2020, plus a handful of newly-populated cells) ---------------------
$updates = @{
    "M2"   = 7
    "B3"   = 2
    "M3"   = 2
    "B5"   = 10
    "X5"   = 16
    "M6"   = 6
    "X7"   = 6
    "B8"   = 3
    "B9"   = 2
    "X10"  = 5
    "X11"  = 19
    "AI11" = 6
    "BE11" = 5
    "BP11" = 2
    "M13"  = 2
    "M15"  = 4
    "B20"  = 2
    "M22"  = 12
    "X26"  = 8
    "X31"  = 2
    "M32"  = 2
    "AT32" = 1
    "X36"  = 2
    "BP41" = 5
    "BP44" = 2
    "B53"  = 1
    "X53"  = 2
    "X67"  = 1
    "B74"  = 1
    "AT76" = 1
    "CA93" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
